$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 139 (shifts old rows 139-163 down to 141-165)
$ws.Rows.Item(139).Insert()
$ws.Rows.Item(139).Insert()

# New row 139: Zafiro rojo
$ws.Cells.Item(139, 1).Value = 7
$ws.Cells.Item(139, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(139, 3).Value = "Ñuble"
$ws.Cells.Item(139, 4).Value = 44491
$ws.Cells.Item(139, 5).Value = 16
$ws.Cells.Item(139, 6).Value = 100112002
$ws.Cells.Item(139, 7).Value = "Pimiento"
$ws.Cells.Item(139, 8).Value = "Zafiro rojo"
$ws.Cells.Item(139, 9).Value = "Primera"
$ws.Cells.Item(139, 10).Value = 100
$ws.Cells.Item(139, 11).Value = 43000
$ws.Cells.Item(139, 12).Value = 44000
$ws.Cells.Item(139, 13).Value = 43500
$ws.Cells.Item(139, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(139, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(139, 16).Value = 2900
$ws.Cells.Item(139, 17).Value = 15
$ws.Cells.Item(139, 18).Value = "Hortaliza"

# New row 140: Zafiro verde
$ws.Cells.Item(140, 1).Value = 7
$ws.Cells.Item(140, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(140, 3).Value = "Ñuble"
$ws.Cells.Item(140, 4).Value = 44491
$ws.Cells.Item(140, 5).Value = 16
$ws.Cells.Item(140, 6).Value = 100112002
$ws.Cells.Item(140, 7).Value = "Pimiento"
$ws.Cells.Item(140, 8).Value = "Zafiro verde"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 100
$ws.Cells.Item(140, 11).Value = 41000
$ws.Cells.Item(140, 12).Value = 42000
$ws.Cells.Item(140, 13).Value = 41500
$ws.Cells.Item(140, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(140, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(140, 16).Value = 2767
$ws.Cells.Item(140, 17).Value = 15
$ws.Cells.Item(140, 18).Value = "Hortaliza"
